$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" date field from 12/19/2025 to
#    12/21/2025 everywhere it is rendered: the slide master and every
#    slide layout's Date placeholder.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "12/19/2025") {
                $shp.TextFrame.TextRange.Text = "12/21/2025"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Fix the swapped c/h axis labels on the OKLCh diagram (slide 1).
#    The big diagram's two standalone labels and the small inset diagram's
#    two labels (inside "Group 147") had their "c" and "h" text swapped.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$bigC = $slide1.Shapes.Item("TextBox 122")
$bigH = $slide1.Shapes.Item("TextBox 124")
Write-Host "Big diagram labels before: $($bigC.Name)='$($bigC.TextFrame.TextRange.Text)' $($bigH.Name)='$($bigH.TextFrame.TextRange.Text)'"
$bigC.TextFrame.TextRange.Text = "h"
$bigH.TextFrame.TextRange.Text = "c"

$group = $slide1.Shapes.Item("Group 147")
$smallC = $group.GroupItems.Item("TextBox 142")
$smallH = $group.GroupItems.Item("TextBox 143")
Write-Host "Inset diagram labels before: $($smallC.Name)='$($smallC.TextFrame.TextRange.Text)' $($smallH.Name)='$($smallH.TextFrame.TextRange.Text)'"
$smallC.TextFrame.TextRange.Text = "h"
$smallH.TextFrame.TextRange.Text = "c"

Write-Host "Done: dates bumped to 12/21/2025 and c/h labels fixed on OKLCh diagram."
